# Commit: "added 4wk low sales check"
#
# The MyForecast model was re-run with a 4-week-low-sales guard, which
# shifted the weekly forecast (column D, "MyForecast") and the derived
# Seasonality Index (column L) on the "Forecast Comparison" sheet. The
# "Summary" sheet's forecast-total rollups are refreshed to match.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison: MyForecast (D) + Seasonality Index (L) ---
# These are numeric cells in both the before and after workbook, so a
# plain numeric .Value assignment keeps their cell type intact.

$wsForecast.Range("D2").Value  = 39
$wsForecast.Range("L2").Value  = 1.18

$wsForecast.Range("D3").Value  = 39
$wsForecast.Range("L3").Value  = 1.19

$wsForecast.Range("D4").Value  = 38
$wsForecast.Range("L4").Value  = 1.18

$wsForecast.Range("D5").Value  = 37
$wsForecast.Range("L5").Value  = 1.04

$wsForecast.Range("D6").Value  = 37
$wsForecast.Range("L6").Value  = 0.83

$wsForecast.Range("D7").Value  = 37
$wsForecast.Range("L7").Value  = 0.91

$wsForecast.Range("D8").Value  = 36
$wsForecast.Range("L8").Value  = 1.12

$wsForecast.Range("D9").Value  = 36
$wsForecast.Range("L9").Value  = 0.85

$wsForecast.Range("D10").Value = 35
$wsForecast.Range("L10").Value = 0.9399999999999999

$wsForecast.Range("D11").Value = 35
$wsForecast.Range("L11").Value = 1

$wsForecast.Range("D12").Value = 35
$wsForecast.Range("L12").Value = 0.84

$wsForecast.Range("L13").Value = 0.93

$wsForecast.Range("D14").Value = 34
$wsForecast.Range("L14").Value = 1.2

$wsForecast.Range("D15").Value = 33
$wsForecast.Range("L15").Value = 0.89

$wsForecast.Range("L16").Value = 1.17

$wsForecast.Range("D17").Value = 33
$wsForecast.Range("L17").Value = 0.99

# --- Summary: refreshed forecast totals ---
# These cells are stored as text in the source workbook (e.g. "540", not
# 540), and a direct $range.Value = "<numeric string>" assignment would be
# auto-coerced to a number by Excel's normal type inference. To preserve
# the original text cell type we stage the new text in a scratch cell via
# a TEXT() formula, copy it, and paste *values only* into the destination
# - a paste-values operation carries the text through verbatim instead of
# re-inferring its type - then clear the scratch cell.

$scratch = $wsSummary.Range("Z1")

function Set-AsText($targetAddress, $number) {
    $scratch.Formula = "=TEXT(" + $number + ",""0"")"
    $scratch.Copy()
    $wsSummary.Range($targetAddress).PasteSpecial(-4163)
    $scratch.ClearContents()
}

Set-AsText "B9"  572
Set-AsText "B10" 299
Set-AsText "B11" 153
Set-AsText "B12" 39
Set-AsText "B14" 33
